$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Calendar2021")

# Row 20: append new slides/video links to the "Theorm Proving and Resolution" topic
$ws.Range("D20").Value = "Theorm Proving and Resolution;slides(slides/14_PropLogic_Part2.pdf); video(https://canvas.jmu.edu/courses/1775272/modules)"

# Row 21: append new slides/video links to the "FOL, Unification" topic, and add a reading entry
$ws.Range("D21").Value = "FOL, Unification;slides(slides/15_FOL_Part1.pdf); video(https://canvas.jmu.edu/courses/1775272/modules)"
$ws.Range("E21").Value = "Chp 8.1 - 8.2"

# Row 22: add a reading entry (topic text unchanged)
$ws.Range("E22").Value = "Chp 8.3"

# Row 23: replace "Prob" topic with "Probability Review/Lab"
$ws.Range("D23").Value = "Probability Review/Lab"

# Update the active selection to D24 (matches the saved view state in the workbook)
$ws.Activate() | Out-Null
$ws.Range("D24").Select() | Out-Null
